$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.408.96'
$ws.Range("E2").Value = '  -2.80%  '
$ws.Range("D3").Value = '1.748.64'
$ws.Range("E3").Value = '  -3.73%  '
$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''321.36'
$ws.Range("E5").Value = '  -2.54%  '
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").Value = '''0.4228'
$ws.Range("E7").Value = '  -4.68%  '
$ws.Range("D8").Value = '''0.3594'
$ws.Range("E8").Value = '  -2.82%  '
$ws.Range("D9").Value = '''0.07535'
$ws.Range("E9").Value = '  -2.06%  '
$ws.Range("D10").Value = '''42.38'
$ws.Range("E10").Value = '  -4.92%  '
$ws.Range("D11").Value = '''1.100'
$ws.Range("E11").Value = '  -2.50%  '
$ws.Range("D12").Value = '''1.001'
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("E13").Value = '  -6.83%  '
$ws.Range("D14").Value = '''6.028'
$ws.Range("E14").Value = '  -4.06%  '
$ws.Range("D15").Value = '''7.215'
$ws.Range("E15").Value = '  -4.99%  '
$ws.Range("D16").Value = '1.750.01'
$ws.Range("E16").Value = '  -4.94%  '
$ws.Range("D17").Value = '''91.39'
$ws.Range("D18").Value = '''0.00001069'
$ws.Range("E18").Value = '  -1.40%  '
$ws.Range("D19").Value = '''0.06359'
$ws.Range("E19").Value = '  -3.62%  '
$ws.Range("D20").Value = '''1.001'
$ws.Range("E20").Value = '  +0.07%  '
$ws.Range("D21").Value = '''17.01'
$ws.Range("E21").Value = '  -3.00%  '
$ws.Range("D22").Value = '''5.879'
$ws.Range("E22").Value = '  -5.41%  '
$ws.Range("D23").Value = '27.443.53'
$ws.Range("E23").Value = '  -2.88%  '
$ws.Range("D24").Value = '''11.18'
$ws.Range("E24").Value = '  -4.40%  '
$ws.Range("D25").Value = '''2.089'
$ws.Range("E25").Value = '  -3.24%  '
$ws.Range("D26").Value = '''160.94'
$ws.Range("E26").Value = '  +3.01%  '
$ws.Range("D27").Value = '''20.25'
$ws.Range("E27").Value = '  -2.54%  '
$ws.Range("D28").Value = '1.948.31'
$ws.Range("E28").Value = '  -4.42%  '
$ws.Range("D29").Value = '''2.130'
$ws.Range("E29").Value = '  -8.63%  '
$ws.Range("D30").Value = '''123.33'
$ws.Range("E30").Value = '  -3.89%  '
$ws.Range("D31").Value = '''1.113'
$ws.Range("E31").Value = '  -7.51%  '
$ws.Range("E32").Value = '  -0.38%  '
$ws.Range("D33").Value = '''5.540'
$ws.Range("E33").Value = '  -5.78%  '
$ws.Range("D34").Value = '''0.08863'
$ws.Range("E34").Value = '  -3.92%  '
$ws.Range("D35").Value = '''12.24'
$ws.Range("E35").Value = '  -6.56%  '
$ws.Range("D36").Value = '''0.02274'
$ws.Range("E36").Value = '  -3.43%  '
$ws.Range("D37").Value = '''0.2099'
$ws.Range("E37").Value = '  -3.85%  '
$ws.Range("D38").Value = '''0.06003'
$ws.Range("E38").Value = '  -3.82%  '
$ws.Range("D39").Value = '''0.6330'
$ws.Range("E39").Value = '  -3.72%  '
$ws.Range("D40").Value = '''4.933'
$ws.Range("E40").Value = '  -4.53%  '
$ws.Range("D41").Value = '''1.177'
$ws.Range("E41").Value = '  -1.72%  '
$ws.Range("D42").Value = '''1.000'
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("D43").Value = '''7.880'
$ws.Range("E43").Value = '  -3.57%  '
$ws.Range("D44").Value = '''1.384'
$ws.Range("E44").Value = '  -1.19%  '
$ws.Range("E45").Value = '  -3.32%  '
$ws.Range("D46").Value = '''0.5859'
$ws.Range("E46").Value = '  -3.57%  '
$ws.Range("D47").Value = '''3.684'
$ws.Range("E47").Value = '  -2.26%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '''1.973'
$ws.Range("E48").Value = '  -3.32%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = '''122.30'
$ws.Range("E49").Value = '  -3.85%  '
$ws.Range("D50").Value = '''1.163'
$ws.Range("E50").Value = '  +0.87%  '
$ws.Range("D51").Value = '''0.06797'
$ws.Range("E51").Value = '  -2.60%  '
